# Insert a new data row at row 233 ("Femacal de La Calera" / Apio sheet).
# Excel's Rows.Insert() shifts row 233..274 down to 234..275 and extends
# the used range (dimension) automatically, matching the diff which turns
# A1:R274 into A1:R275.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("233:233").Insert()

$ws.Range("A233").Value = 3
$ws.Range("B233").Value = "Femacal de La Calera"
$ws.Range("C233").Value = "Coquimbo"
$ws.Range("D233").Value = 44522
$ws.Range("E233").Value = 5
$ws.Range("F233").Value = 100112017
$ws.Range("G233").Value = "Apio"
$ws.Range("H233").Value = "Americana (o)"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 230
$ws.Range("K233").Value = 8500
$ws.Range("L233").Value = 9000
$ws.Range("M233").Value = 8826
$ws.Range("N233").Value = "`$/docena de matas"
$ws.Range("O233").Value = "Pan de Az$([char]0xFA)car"
$ws.Range("P233").Value = 1471
$ws.Range("Q233").Value = 6
$ws.Range("R233").Value = "Hortaliza"
